$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.302.92"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "1.859.74"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.19%  "
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4757"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.69%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2767"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06446"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "1.867.78"
$ws.Range("E10").Value = "  -0.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07439"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("E12").Value = "  -3.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.011"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "85.51"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6335"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.89%  "
$ws.Range("D16").Value = "30.273.56"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007340"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.58%  "
$ws.Range("D20").Value = "2.110.03"
$ws.Range("E20").Value = "  -0.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "224.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.093"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.034"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.226"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.868"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1033"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.376"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.218"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.902"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04903"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.16%  "
$ws.Range("E34").Value = "  -4.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7285"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.05%  "
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.715"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01902"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.624"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9029"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.984"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9956"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4102"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.24%  "
$ws.Range("E45").Value = "  -5.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.067"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "61.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1208"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.808"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.399"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05603"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.48%  "
